# AUTOSAR_SWS_CANInterface.xlsx - template update for req/design/test case
# Populate Sheet1 with the requirement-tracking header row, formatting,
# column widths, frozen header pane and per-column explanatory comments.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row values (land in a shared string table) -------------------
$ws.Range("A1").Value = "ReqID"
$ws.Range("B1").Value = "Requirements"
$ws.Range("C1").Value = "Category"
$ws.Range("D1").Value = "Status"
$ws.Range("E1").Value = "TestCase ID"
$ws.Range("F1").Value = "Design ID"
$ws.Range("G1").Value = "Review Status"

# --- Header row formatting: bold font on a shaded fill --------------------
$header = $ws.Range("A1:G1")
$header.Font.Bold = $true
$header.Interior.ThemeColor = 4
$header.Interior.TintAndShade = -0.099978637043366805

# --- Column widths ---------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 68.16666666666667
$ws.Columns.Item(3).ColumnWidth = 10.833333333333332
$ws.Columns.Item(5).ColumnWidth = 12
$ws.Columns.Item(6).ColumnWidth = 10
$ws.Columns.Item(7).ColumnWidth = 11.5

# --- Freeze the header row and leave the selection on B5 -------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B5").Select()

# --- Printing ---------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# --- Explanatory comments on the header cells -------------------------------
$excel.UserName = "Author"

$c = $ws.Range("A1").AddComment()
$c.Text("Author:" + [char]10 + "Requirement Identifier unique for each new line")

$c = $ws.Range("C1").AddComment()
$c.Text("Author:" + [char]10 + "Type of requirement functional or non functional")

$c = $ws.Range("D1").AddComment()
$c.Text("Author:" + [char]10 + "Status of implementation" + [char]10 + "OPEN, NA, CLOSED")

$c = $ws.Range("E1").AddComment()
$c.Text("Author:" + [char]10 + "Mention the linked Test Case ID and document name with version number")

$c = $ws.Range("F1").AddComment()
$c.Text("Author:" + [char]10 + "Mention the linked Design ID and document name with version number")

$c = $ws.Range("G1").AddComment()
$c.Text("Author:" + [char]10 + "Review status , INREVIEW, COMPLETED, FORREVIEW, OPEN")

# --- Sheet2 / Sheet3 keep their default row height (template placeholders) -
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")
